{"js": "// Apply the README-stats / docx-prep fix for the\n// Renaissance / JDK17 / ZGC / philosophers heap-8G results table.\n//\n// Summary of the change:\n//   - Row 1 (\"100\")  -> \"0M\"\n//   - Row 2 (\"0.01\") -> \"0M\"\n//   - Row 3 (\"181\")  -> \"0M\"\n//   - 10 brand-new rows inserted right after (old) row 3, carrying the\n//     values that used to be crammed (tab-separated) into the final three\n//     rows of the table: 189, 0.00002, 0.00005, 0.00003, 0.00001, 0.00003,\n//     0.00003, 0.00004, 0.00519, 100.0\n//   - The last three rows of the table, which held 10 tab-separated values\n//     each in a single run, are collapsed down to a single plain value:\n//       \"86\\t0.00002\\t...\\t100.0\"  -> \"100\"\n//       \"100\\t0.00002\\t...\\t100.0\" -> \"0.01\"\n//       \"3\\t0.00002\\t...\\t100.0\"   -> \"181\"\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst rowCountBefore = rows.items.length;\n\n// --- Step 1: fix up the first three rows -------------------------------\nrows.items[0].cells.load(\"items\");\nrows.items[1].cells.load(\"items\");\nrows.items[2].cells.load(\"items\");\nawait context.sync();\n\nrows.items[0].cells.items[0].value = \"0M\";\nrows.items[1].cells.items[0].value = \"0M\";\nrows.items[2].cells.items[0].value = \"0M\";\nawait context.sync();\n\n// --- Step 2: insert the 10 new rows right after (old) row index 2 ------\nconst newValues = [\n  [\"189\"],\n  [\"0.00002\"],\n  [\"0.00005\"],\n  [\"0.00003\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00519\"],\n  [\"100.0\"],\n];\nrows.items[2].insertRows(\"After\", newValues.length, newValues);\nawait context.sync();\n\n// --- Step 3: collapse the last three (tab-separated) rows down to a ----\n// single value each. Re-load rows since the collection changed size.\nrows.load(\"items\");\nawait context.sync();\n\nconst lastCollapsed = [\"100\", \"0.01\", \"181\"];\nconst total = rows.items.length; // rowCountBefore + newValues.length\nconst lastThreeStart = total - 3;\n\nfor (let i = 0; i < 3; i++) {\n  rows.items[lastThreeStart + i].cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let i = 0; i < 3; i++) {\n  rows.items[lastThreeStart + i].cells.items[0].value = lastCollapsed[i];\n}\nawait context.sync();\n", "ps1": "# Apply the README-stats / docx-prep fix for the\n# Renaissance / JDK17 / ZGC / philosophers heap-8G results table.\n#\n# Summary of the change:\n#   - Row 1 (\"100\")  -> \"0M\"\n#   - Row 2 (\"0.01\") -> \"0M\"\n#   - Row 3 (\"181\")  -> \"0M\"\n#   - 10 brand-new rows inserted right after (old) row 3, carrying the\n#     values that used to be crammed (tab-separated) into the final three\n#     rows of the table: 189, 0.00002, 0.00005, 0.00003, 0.00001, 0.00003,\n#     0.00003, 0.00004, 0.00519, 100.0\n#   - The last three rows of the table, which held 10 tab-separated values\n#     each in a single run, are collapsed down to a single plain value:\n#       \"86\\t0.00002\\t...\\t100.0\"  -> \"100\"\n#       \"100\\t0.00002\\t...\\t100.0\" -> \"0.01\"\n#       \"3\\t0.00002\\t...\\t100.0\"   -> \"181\"\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Step 1: fix up the first three rows --------------------------------\n$t.Rows.Item(1).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(2).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(3).Cells.Item(1).Range.Text = \"0M\"\n\n# --- Step 2: insert 10 new rows right after (old) row 3 -----------------\n# Table.Rows.Add(beforeRow) always inserts immediately BEFORE beforeRow,\n# so walk the new values back-to-front, always inserting before the same\n# anchor (the row that used to be row 4) to land them in the right order.\n$newValues = @(\"189\", \"0.00002\", \"0.00005\", \"0.00003\", \"0.00001\", \"0.00003\", \"0.00003\", \"0.00004\", \"0.00519\", \"100.0\")\n$anchorRow = $t.Rows.Item(4)\nfor ($i = $newValues.Length - 1; $i -ge 0; $i--) {\n  $newRow = $t.Rows.Add($anchorRow)\n  $newRow.Cells.Item(1).Range.Text = $newValues[$i]\n}\n\n# --- Step 3: collapse the last three (tab-separated) rows down to a -----\n# single value each.\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = \"100\"\n$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = \"0.01\"\n$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = \"181\"\n"}
